$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1823923333333334
$ws.Range("H2").Value = 0.547177
$ws.Range("M2").Value = 2.133443333333334
$ws.Range("N2").Value = 6.40033
$ws.Range("O2").Value = 0.2605947899689859
$ws.Range("P2").Value = 0.2605947899689859
$ws.Range("Q2").Value = 0.3891237076011112
$ws.Range("R2").Value = 3.50211336841
$ws.Range("S2").Value = 0.2605947899689859
$ws.Range("T2").Value = 0.2605947899689859

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1823923333333334
$ws.Range("H3").Value = 0.547177
$ws.Range("O3").Value = 0.5209338844846115
$ws.Range("P3").Value = 0.5209338844846116
$ws.Range("Q3").Value = 0.7778656072511112
$ws.Range("R3").Value = 7.000790465260001
$ws.Range("S3").Value = 0.5209338844846115
$ws.Range("T3").Value = 0.5209338844846116

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1823923333333334
$ws.Range("H4").Value = 0.547177
$ws.Range("M4").Value = 1.788586
$ws.Range("N4").Value = 5.365758
$ws.Range("O4").Value = 0.2184713255464024
$ws.Range("P4").Value = 0.2184713255464024
$ws.Range("Q4").Value = 0.3262243739073333
$ws.Range("R4").Value = 2.936019365166
$ws.Range("S4").Value = 0.2184713255464024
$ws.Range("T4").Value = 0.2184713255464024
